$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("A1").Value = "Filiere"
$ws.Range("B1").Value = "Secteur"

# Convert the range into a formatted Excel Table
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:B2"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium1"

# Update selection
[void]$ws.Range("B2").Select()

# Adjust column widths (closest the engine's width grid can represent to the
# authored 22.109375 / 20.88671875 character widths)
$ws.Columns.Item(1).ColumnWidth = 22.0
$ws.Columns.Item(2).ColumnWidth = 20.0
